$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Medium" table: append the 32-node / 1280-process benchmark row ---
# (new row 19, inserted right after the existing 16-node row 18)
$ws.Range("A19").Value = 32
$ws.Range("C19").Value = 1280
$ws.Range("F19").Value = 3233
$ws.Range("J19").Value = 1448

# --- "Large" table: append the matching 32-node / 1280-process row ---
# (new row 38, appended after the existing 16-node row 37)
$ws.Range("A38").Value = 32
$ws.Range("C38").Value = 1280
$ws.Range("F38").Value = 4002
$ws.Range("J38").Value = 2297

# --- Match the workbook's saved view/selection state ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("K38").Select()
